$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = 45659
$ws.Range("A28").NumberFormat = $ws.Range("A27").NumberFormat
$ws.Range("B28").Value = "testing, polishing and rat"
$ws.Range("C28").Value = 5

$ws.Range("F20").Select()
